$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 420
$ws1.Range("F5").Value = 461
$ws1.Range("F7").Value = 2491
$ws1.Range("F9").Value = 6574
$ws1.Range("F10").Value = 176
$ws1.Range("F11").Value = 422
$ws1.Range("F12").Value = 29

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 420
$ws4.Range("F5").Value = 461
$ws4.Range("F9").Value = 2491
$ws4.Range("F11").Value = 6574
$ws4.Range("F12").Value = 176
$ws4.Range("F13").Value = 422
$ws4.Range("F16").Value = 29
